$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-10 Sunday" "2024-03-11 Monday"

Replace-Text "36×30=1080" "62×48=2976"
Replace-Text "21×47=987" "87×34=2958"
Replace-Text "14×15=210" "66×20=1320"
Replace-Text "81×53=4293" "76×38=2888"
Replace-Text "49×38=1862" "51×62=3162"

Replace-Text "99×60=5940" "33×26=858"
Replace-Text "88×12=1056" "13×83=1079"
Replace-Text "28×57=1596" "38×90=3420"
Replace-Text "39×88=3432" "53×43=2279"
Replace-Text "37×53=1961" "37×18=666"

Replace-Text "71×96=6816" "41×83=3403"
Replace-Text "69×96=6624" "55×82=4510"
Replace-Text "63×22=1386" "98×65=6370"
Replace-Text "71×93=6603" "47×84=3948"
Replace-Text "89×36=3204" "91×30=2730"

Replace-Text "39×76=2964" "68×95=6460"
Replace-Text "70×70=4900" "17×88=1496"
Replace-Text "18×82=1476" "72×40=2880"
Replace-Text "36×86=3096" "91×66=6006"
Replace-Text "83×78=6474" "64×87=5568"

Replace-Text "49×54=2646" "51×73=3723"
Replace-Text "12×80=960" "11×72=792"
Replace-Text "76×31=2356" "44×68=2992"
Replace-Text "66×22=1452" "68×31=2108"
Replace-Text "80×67=5360" "46×51=2346"
